# Update Name of Algo
# Applies updated imputation results to specific cells in Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -22.13250000000001
$ws.Range("D10").Value = -7.889400000000003
$ws.Range("D15").Value = -8.314899999999996
$ws.Range("C18").Value = -11.47439999999999
$ws.Range("A21").Value = -19.77039999999999
$ws.Range("D21").Value = -7.840299999999998
$ws.Range("D22").Value = -8.168000000000005
$ws.Range("A23").Value = -20.50799999999998
$ws.Range("B24").Value = 5.653300000000002
$ws.Range("D24").Value = -6.850799999999996
$ws.Range("A25").Value = -21.871
$ws.Range("B28").Value = 6.512800000000001
$ws.Range("B36").Value = 9.459200000000006
$ws.Range("B45").Value = 5.043100000000003
$ws.Range("D46").Value = -7.818899999999996
$ws.Range("B48").Value = 5.653700000000004
$ws.Range("B49").Value = 5.866899999999996
$ws.Range("C51").Value = -12.17920000000001
$ws.Range("B52").Value = 5.595399999999997
$ws.Range("A53").Value = -21.94289999999999
$ws.Range("B53").Value = 5.821700000000002
$ws.Range("B54").Value = 4.999400000000003
$ws.Range("C55").Value = -13.11249999999999
$ws.Range("D56").Value = -8.697400000000002
$ws.Range("A57").Value = -22.2009
$ws.Range("A59").Value = -22.20790000000001
$ws.Range("D61").Value = -8.243199999999998
$ws.Range("C64").Value = -10.60319999999999
$ws.Range("D66").Value = -6.983199999999997
$ws.Range("A69").Value = -21.65539999999998
$ws.Range("B70").Value = 6.888100000000004
$ws.Range("D74").Value = -8.339400000000007
$ws.Range("A79").Value = -20.08110000000001
$ws.Range("C80").Value = -13.11360000000001
$ws.Range("A83").Value = -21.895
$ws.Range("B86").Value = 5.1127
$ws.Range("B87").Value = 5.315499999999997
$ws.Range("D87").Value = -8.392399999999995
$ws.Range("D88").Value = -7.980199999999991
$ws.Range("C92").Value = -10.50740000000001
$ws.Range("A93").Value = -21.20660000000002
$ws.Range("C94").Value = -10.3815
$ws.Range("C96").Value = -10.16940000000001
$ws.Range("D100").Value = -8.195100000000002
$ws.Range("B101").Value = 5.1887
